$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------
# 1. Re-style existing row 9 (B9:F9) to the "highlighted" look
#    (same visual treatment as row 7) by copying formats from
#    cells that already carry the target cellXfs indices.
#    B9 -> s=10, C9 -> s=11, D9 -> s=10, E9/F9 -> s=13
# ---------------------------------------------------------------
$ws2.Range("B7").Copy()
$ws2.Range("B9").PasteSpecial(-4122)

$ws2.Range("C7").Copy()
$ws2.Range("C9").PasteSpecial(-4122)

$ws2.Range("B7").Copy()
$ws2.Range("D9").PasteSpecial(-4122)

$ws1.Range("F9").Copy()
$ws2.Range("E9").PasteSpecial(-4122)
$ws2.Range("F9").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# ---------------------------------------------------------------
# 2. Add the new row (row 10) : "Painter's partition problem"
# ---------------------------------------------------------------
$ws2.Range("B10").Value = 6
$ws2.Range("C10").Value = "Binary S2 9"
$ws2.Range("D10").Value = "Painter's partition problem"

# ---------------------------------------------------------------
# 3. Hyperlink for E10 (must be created before the final formats
#    are applied, because adding a hyperlink / setting its value
#    resets the cell's number format / style)
# ---------------------------------------------------------------
$url = "https://www.scaler.com/academy/mentee-dashboard/class/30366/assignment/problems/271?navref=cl_tt_lst_nm"
$link = $ws2.Hyperlinks.Add($ws2.Range("E10"), $url)
$link.TextToDisplay = $url
$ws2.Range("E10").Value = "Painter's Partition Problem - Problem | Scaler Academy"

$ws2.Range("G10").Value = "Special modular division"

# Apply matching cell formats (reuse existing style definitions)
$ws2.Range("B1").Copy()
$ws2.Range("B10").PasteSpecial(-4122)

$ws2.Range("B1").Copy()
$ws2.Range("C10").PasteSpecial(-4122)

$ws2.Range("D3").Copy()
$ws2.Range("D10").PasteSpecial(-4122)

$ws1.Range("E11").Copy()
$ws2.Range("E10").PasteSpecial(-4122)

$ws2.Range("D3").Copy()
$ws2.Range("G10").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Row height for the new row (matches the other wrapped-text rows)
$ws2.Rows.Item(10).RowHeight = 43.2

# ---------------------------------------------------------------
# 4. Update selection / active cell to reflect the new last row
# ---------------------------------------------------------------
$ws2.Activate()
$ws2.Range("C10").Select()
